$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "GRT-USD"
$ws.Range("A15").Value = "BSCX-USD"
